$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.011.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.650.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3928"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3886"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.363"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.008"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08480"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.202"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.886"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001314"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.649.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.915"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.009"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.994.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.518"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.068"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.862"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.831.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.038"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03024"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.728"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2721"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09171"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.432"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6940"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.486"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.096"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08293"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.408"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.08%  "
